# MassWateR ExampleSites.xlsx update:
#  1. Insert a new "ABT-162 / Cox Street bridge" site row into the Sites table
#     (alphabetically between ABT-144 and ABT-237).
#  2. Add a new "Instructions" worksheet (placed after "Sites") that documents
#     the Sites-tab field requirements.
#  3. Re-point the selection / active cell bookkeeping to match.

$wb = $excel.ActiveWorkbook
$sites = $wb.Worksheets.Item("Sites")

# ---------------------------------------------------------------------------
# 1. Sites sheet: insert the new ABT-162 row (row 6, pushing ABT-237.. down)
# ---------------------------------------------------------------------------
$sites.Rows.Item(6).Insert()
$sites.Range("A6").Value = "ABT-162"
$sites.Range("B6").Value = "Cox Street bridge"
$sites.Range("C6").Value = 42.399797
$sites.Range("D6").Value = -71.545985000000002
$sites.Range("E6").Value = "Assabet"

# ---------------------------------------------------------------------------
# 2. Add the Instructions worksheet right after Sites
# ---------------------------------------------------------------------------
$instr = $wb.Worksheets.Add($null, $sites)
$instr.Name = "Instructions"

# Column widths (character units)
$instr.Columns.Item(1).ColumnWidth = 30.28515625
$instr.Columns.Item(2).ColumnWidth = 96.42578125
$instr.Columns.Item(3).ColumnWidth = 21.42578125
$instr.Columns.Item(4).ColumnWidth = 21.5703125
$instr.Columns.Item(5).ColumnWidth = 25

# --- Row 1-2: callouts -------------------------------------------------
$instr.Range("A1").Value = "The Sites tab must be formatted exactly like the Sites template, with all of the following fields."
$instr.Range("A1").Font.Bold = $true
$instr.Range("A1").Font.Color = 0xC07000

$instr.Range("C1").Value = "Template updated 5/19/23"
$instr.Range("C1").Font.Color = 0x0000FF

$instr.Range("A2").Value = "The Sites tab must be the first tab in this workbook."
$instr.Range("A2").Font.Bold = $true
$instr.Range("A2").Font.Color = 0xC07000

$instr.Range("C2").Value = "Samples updated 5/19/23"
$instr.Range("C2").Font.Color = 0x0000FF

# --- Row 4: table header --------------------------------------------------
$headerValues = @("Field", "Instructions", "Example", "Available Values", "Required?")
for ($i = 0; $i -lt 5; $i++) {
    $col = [char](65 + $i)
    $cell = $instr.Range("$col4")
    $cell.Value = $headerValues[$i]
    $cell.HorizontalAlignment = -4108
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(7).Weight = 2
    $cell.Borders.Item(10).LineStyle = 1
    $cell.Borders.Item(10).Weight = 2
    $cell.Borders.Item(8).LineStyle = 1
    $cell.Borders.Item(8).Weight = 2
    $cell.Borders.Item(9).LineStyle = 1
    $cell.Borders.Item(9).Weight = -4138
}
$instr.Rows.Item(4).RowHeight = 15.75

# --- Rows 5-9: field reference table --------------------------------------
$fieldRows = @(
    @{ Row=5; A="Monitoring Location ID"; B="Location ID that is used in your Results file.  Must match exactly."; C="ABT-010"; D="any"; E="Required" },
    @{ Row=6; A="Monitoring Location Name"; B="Name of monitoring location."; C="477 Lowell Rd, Concord"; D="any"; E="Required for WQX" },
    @{ Row=7; A="Monitoring Location Latitude "; B="Latitude of monitoring location in decimal form.  At least 5 decimals."; C=42.470370000000003; D="any"; E="Required for mapping" },
    @{ Row=8; A="Monitoring Location Longitude"; B="Longitude of monitoring location in decimal form.  At least 5 decimals."; C=-71.362578999999997; D="any"; E="Required for mapping" },
    @{ Row=9; A="Location Group"; B="An optional free-form grouping attribute.  This will allow you to summarize locations by group in the graphing and mapping analysis functions."; C="Lower Assabet"; D="any"; E="Optional" }
)

foreach ($fr in $fieldRows) {
    $r = $fr.Row

    $a = $instr.Range("A$r")
    $a.Value = $fr.A
    $a.VerticalAlignment = -4160
    $a.Borders.Item(7).LineStyle = 1
    $a.Borders.Item(7).Weight = 2
    $a.Borders.Item(10).LineStyle = 1
    $a.Borders.Item(10).Weight = 2
    $a.Borders.Item(9).LineStyle = 1
    $a.Borders.Item(9).Weight = 2
    if ($r -ne 5) {
        $a.Borders.Item(8).LineStyle = 1
        $a.Borders.Item(8).Weight = 2
    }

    $b = $instr.Range("B$r")
    $b.Value = $fr.B
    $b.WrapText = $true
    $b.VerticalAlignment = -4160
    $b.Borders.Item(7).LineStyle = 1
    $b.Borders.Item(7).Weight = 2
    $b.Borders.Item(10).LineStyle = 1
    $b.Borders.Item(10).Weight = 2
    $b.Borders.Item(8).LineStyle = 1
    $b.Borders.Item(8).Weight = 2
    $b.Borders.Item(9).LineStyle = 1
    $b.Borders.Item(9).Weight = 2

    $c = $instr.Range("C$r")
    $c.Value = $fr.C
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4160
    $c.Borders.Item(7).LineStyle = 1
    $c.Borders.Item(7).Weight = 2
    $c.Borders.Item(10).LineStyle = 1
    $c.Borders.Item(10).Weight = 2
    $c.Borders.Item(8).LineStyle = 1
    $c.Borders.Item(8).Weight = 2
    $c.Borders.Item(9).LineStyle = 1
    $c.Borders.Item(9).Weight = 2

    $d = $instr.Range("D$r")
    $d.Value = $fr.D
    $d.HorizontalAlignment = -4108
    $d.VerticalAlignment = -4160
    $d.Font.Italic = $true
    $d.Borders.Item(7).LineStyle = 1
    $d.Borders.Item(7).Weight = 2
    $d.Borders.Item(10).LineStyle = 1
    $d.Borders.Item(10).Weight = 2
    $d.Borders.Item(8).LineStyle = 1
    $d.Borders.Item(8).Weight = 2
    $d.Borders.Item(9).LineStyle = 1
    $d.Borders.Item(9).Weight = 2

    $e = $instr.Range("E$r")
    $e.Value = $fr.E
    $e.HorizontalAlignment = -4108
    $e.VerticalAlignment = -4160
    $e.Font.Italic = $true
    $e.Borders.Item(7).LineStyle = 1
    $e.Borders.Item(7).Weight = 2
    $e.Borders.Item(10).LineStyle = 1
    $e.Borders.Item(10).Weight = 2
    $e.Borders.Item(8).LineStyle = 1
    $e.Borders.Item(8).Weight = 2
    $e.Borders.Item(9).LineStyle = 1
    $e.Borders.Item(9).Weight = 2
}

$instr.Rows.Item(9).RowHeight = 30

# --- Row 11: trailing bold spacer cell ------------------------------------
$instr.Range("B11").Font.Bold = $true

# ---------------------------------------------------------------------------
# 3. Freeze panes on the Instructions sheet & restore selections
# ---------------------------------------------------------------------------
$instr.Activate()
$instr.Range("B5").Select()
$excel.ActiveWindow.FreezePanes = $true
$instr.Range("C3").Select()

$sites.Activate()
$sites.Range("B15").Select()
